$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A144 text
$ws.Range("A144").Value = "Watching."

# Insert a new row at 165, shifting old rows 165-182 down to 166-183.
# This splits the old combined sentence in row 164 into two rows (164 & 165),
# and carries a couple of small text/timestamp corrections further down.
$ws.Rows.Item(165).Insert()

$ws.Range("A164").Value = "Can you fly back?"
$ws.Range("B164").Value = "00:13:46"
$ws.Range("C164").Value = "00:13:46"
$ws.Range("A165").Value = "He fly up?"
$ws.Range("B165").Value = "00:13:46"
$ws.Range("C165").Value = "00:13:47"
$ws.Range("A166").Value = "No, it's."
$ws.Range("B166").Value = "00:13:47"
$ws.Range("C166").Value = "00:13:47"
$ws.Range("A167").Value = "We had a pretty rough landing due to, of all things, pilot error."
$ws.Range("B167").Value = "00:13:49"
$ws.Range("C167").Value = "00:13:52"
$ws.Range("A168").Value = "The plan was to free fall until right about 200 meters above the ground, max out the thrust, and fly back towards the launch site."
$ws.Range("B168").Value = "00:13:52"
$ws.Range("C168").Value = "00:13:59"
$ws.Range("A169").Value = "But due to the stress, and because the landscape looked all the same, I didn't realize that the ground was coming until way too late."
$ws.Range("B169").Value = "00:13:59"
$ws.Range("C169").Value = "00:14:06"
$ws.Range("A170").Value = "Luckily, we knew the general location of the Cansat thanks to the onboard footage, and it only took us a mere three to four minutes before it was located."
$ws.Range("B170").Value = "00:14:06"
$ws.Range("C170").Value = "00:14:15"
$ws.Range("A171").Value = "Somehow, the Kansai was still completely intact with only slight bending in one of the arms."
$ws.Range("B171").Value = "00:14:28"
$ws.Range("C171").Value = "00:14:34"
$ws.Range("A172").Value = "So after getting the arms replaced, the satellite was as good as new."
$ws.Range("B172").Value = "00:14:34"
$ws.Range("C172").Value = "00:14:38"
$ws.Range("A173").Value = "We even managed to recover our data thanks to this single G force indicator."
$ws.Range("B173").Value = "00:14:38"
$ws.Range("C173").Value = "00:14:43"
$ws.Range("A174").Value = "Just make a program to grab the value from every frame, graph it over time, and do some integration."
$ws.Range("B174").Value = "00:14:43"
$ws.Range("C174").Value = "00:14:48"
$ws.Range("A175").Value = "And now you have altitude."
$ws.Range("B175").Value = "00:14:48"
$ws.Range("C175").Value = "00:14:50"
$ws.Range("A176").Value = "We spent the rest of the day just flying this thing around randomly for fun."
$ws.Range("B176").Value = "00:14:50"
$ws.Range("C176").Value = "00:14:54"
$ws.Range("A177").Value = "And we also did this freefall test, and I think I'm now a much better pilot."
$ws.Range("B177").Value = "00:14:54"
$ws.Range("C177").Value = "00:14:59"
$ws.Range("A178").Value = "Speaking of that, we were actually offered a second launch on a bigger and more powerful rocket in Seattle during the summer."
$ws.Range("B178").Value = "00:14:59"
$ws.Range("C178").Value = "00:15:05"
$ws.Range("A179").Value = "So make sure to let me know in the comments if you are interested."
$ws.Range("B179").Value = "00:15:05"
$ws.Range("C179").Value = "00:15:09"
$ws.Range("A180").Value = "And with that being said, special thanks to Zach, Daniel, Emre, Eric and Cynthia for working on this project with me."
$ws.Range("B180").Value = "00:15:13"
$ws.Range("C180").Value = "00:15:19"
$ws.Range("A181").Value = "Special thanks to everyone at Kansat who made this experience possible."
$ws.Range("B181").Value = "00:15:19"
$ws.Range("C181").Value = "00:15:23"
$ws.Range("A182").Value = "And thank you to you for watching."
$ws.Range("B182").Value = "00:15:23"
$ws.Range("C182").Value = "00:15:25"
$ws.Range("A183").Value = "See you next time."
$ws.Range("B183").Value = "00:15:25"
$ws.Range("C183").Value = "00:15:27"
